$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-12-22"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 12-22)"

# Update December (row 13) total-column value and the grand Total row (row 14)
$ws.Range("I13").Value = 98
$ws.Range("I14").Value = 1615
